# Nalco aluminium ingot price sheet: daily refresh.
# A new day's price row is published at the top (row 2); every existing
# row shifts down by one, and the oldest row is re-published once more
# at the bottom (row 115) since no new circular superseded it yet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift all existing data rows (2..114) down by one row; this also
# extends the sheet to row 115 and carries the hyperlinks along with
# their rows.
$ws.Rows.Item(2).Insert()

# New row 2: same circular/price info as the (now shifted) row below it,
# just a day later on the "Date" column.
$ws.Range("A2:A2").NumberFormat = "@"
$ws.Range("A2").Value = "28-11-2025"
$ws.Range("A2:A2").NumberFormat = "General"

$ws.Range("B2").Value = "ALUMINIUM INGOT"
$ws.Range("C2").Value = "IE07"
$ws.Range("D2").Value = 297.15

$ws.Range("E2:E2").NumberFormat = "@"
$ws.Range("E2").Value = "01-11-2025"
$ws.Range("E2:E2").NumberFormat = "General"

$ws.Range("F2").Value = "https://nalcoindia.com/wp-content/uploads/2019/01/Ingot-01-11-2025.pdf"
